$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 37037332
$ws.Range("I33").Value = 99.23529000000001
$ws.Range("J33").Value = 100000620
$ws.Range("K33").Value = 99.23529000000001
$ws.Range("L33").Value = 100000620
$ws.Range("M33").Value = 129.76471
$ws.Range("N33").Value = -100001078

$ws.Range("H41").Value = 394.5
$ws.Range("I41").Value = 232.25
$ws.Range("J41").Value = 502.66666
$ws.Range("K41").Value = 232.25
$ws.Range("L41").Value = 502.66666
$ws.Range("M41").Value = 207.75
$ws.Range("N41").Value = -1382.66666

$ws.Range("H51").Value = 4908.6665
$ws.Range("J51").Value = 4909.4546
$ws.Range("L51").Value = 4909.4546
$ws.Range("N51").Value = -5877.4546

$ws.Range("H106").Value = 38450
$ws.Range("I106").Value = 41865
$ws.Range("K106").Value = 41865
$ws.Range("M106").Value = -41234

$ws.Range("H111").Value = 1258.2667
$ws.Range("I111").Value = 1026.75
$ws.Range("J111").Value = 1522.8572
$ws.Range("K111").Value = 3080.25
$ws.Range("L111").Value = 4568.571599999999
$ws.Range("M111").Value = -13.25
$ws.Range("N111").Value = -10702.5716

$ws.Range("H138").Value = 2452.0984
$ws.Range("I138").Value = 1546.8422
$ws.Range("J138").Value = 3947.739
$ws.Range("K138").Value = 4640.5266
$ws.Range("L138").Value = 11843.217
$ws.Range("M138").Value = 499.4733999999999
$ws.Range("N138").Value = -22123.217

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4355.5
$ws.Range("I61").Value = 1211
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 1211
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -999
$ws.Range("N61").Value = -7924

$ws.Range("H74").Value = 4208.4443
$ws.Range("I74").Value = 3937.652
$ws.Range("K74").Value = 3937.652
$ws.Range("M74").Value = -3063.652

$ws.Range("H77").Value = 4208.4443
$ws.Range("I77").Value = 3937.652
$ws.Range("K77").Value = 19688.26
$ws.Range("M77").Value = -15320.26

$ws.Range("H132").Value = 10135
$ws.Range("I132").Value = 5498.485
$ws.Range("K132").Value = 16495.455
$ws.Range("M132").Value = -13965.455

$ws.Range("H136").Value = 4355.5
$ws.Range("I136").Value = 1211
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 3633
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -1083
$ws.Range("N136").Value = -27600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = ""
$ws.Range("N53").Value = 0

$ws.Range("H94").Value = 1201.3024
$ws.Range("I94").Value = 878.5454999999999
$ws.Range("J94").Value = 1539.4286
$ws.Range("K94").Value = 878.5454999999999
$ws.Range("L94").Value = 1539.4286
$ws.Range("M94").Value = -427.5454999999999
$ws.Range("N94").Value = -2441.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1338.1
$ws.Range("I31").Value = 1147.6875
$ws.Range("J31").Value = 2099.75
$ws.Range("K31").Value = 1147.6875
$ws.Range("L31").Value = 2099.75
$ws.Range("M31").Value = -852.6875
$ws.Range("N31").Value = -2689.75

$ws.Range("H34").Value = 1338.1
$ws.Range("I34").Value = 1147.6875
$ws.Range("J34").Value = 2099.75
$ws.Range("K34").Value = 1147.6875
$ws.Range("L34").Value = 2099.75
$ws.Range("M34").Value = -945.6875
$ws.Range("N34").Value = -2503.75

$ws.Range("H58").Value = 3555.2856
$ws.Range("I58").Value = 3780.1
$ws.Range("K58").Value = 3780.1
$ws.Range("M58").Value = -3577.1

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = ""
$ws.Range("N98").Value = 0

$ws.Range("H136").Value = 3555.2856
$ws.Range("I136").Value = 3780.1
$ws.Range("K136").Value = 11340.3
$ws.Range("M136").Value = -8790.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 452.25
$ws.Range("I107").Value = 203
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 609
$ws.Range("L107").Value = 3600
$ws.Range("M107").Value = 1311
$ws.Range("N107").Value = -7440

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = ""
$ws.Range("N125").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 915.4
$ws.Range("I3").Value = 2089
$ws.Range("J3").Value = 133
$ws.Range("K3").Value = 2089
$ws.Range("L3").Value = 133
$ws.Range("M3").Value = -1973
$ws.Range("N3").Value = -365

$ws.Range("H11").Value = 288601.16
$ws.Range("J11").Value = 288601.16
$ws.Range("L11").Value = 288601.16
$ws.Range("N11").Value = -288879.16

$ws.Range("H13").Value = 596
$ws.Range("I13").Value = 192.5
$ws.Range("J13").Value = 999.5
$ws.Range("K13").Value = 192.5
$ws.Range("L13").Value = 999.5
$ws.Range("M13").Value = -53.5
$ws.Range("N13").Value = -1277.5

$ws.Range("H14").Value = 2176.25
$ws.Range("I14").Value = 1350
$ws.Range("J14").Value = 3002.5
$ws.Range("K14").Value = 1350
$ws.Range("L14").Value = 3002.5
$ws.Range("M14").Value = -1182
$ws.Range("N14").Value = -3338.5

$ws.Range("H17").Value = 4803
$ws.Range("J17").Value = 4803
$ws.Range("L17").Value = 4803
$ws.Range("N17").Value = -5139

$ws.Range("H22").Value = 4399.6
$ws.Range("I22").Value = 6332
$ws.Range("J22").Value = 3571.4285
$ws.Range("K22").Value = 6332
$ws.Range("L22").Value = 3571.4285
$ws.Range("M22").Value = -5803
$ws.Range("N22").Value = -4629.4285

$ws.Range("H23").Value = 1500
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1500
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = ""
$ws.Range("M23").Value = 1500
$ws.Range("N23").Value = -1946

$ws.Range("H25").Value = 1000
$ws.Range("J25").Value = 1000
$ws.Range("L25").Value = 1000
$ws.Range("N25").Value = -2058

$ws.Range("H80").Value = 3829.8462
$ws.Range("I80").Value = 3679.75
$ws.Range("J80").Value = 4070
$ws.Range("K80").Value = 3679.75
$ws.Range("L80").Value = 4070
$ws.Range("M80").Value = -2681.75
$ws.Range("N80").Value = -6066

$ws.Range("H83").Value = 3829.8462
$ws.Range("I83").Value = 3679.75
$ws.Range("J83").Value = 4070
$ws.Range("K83").Value = 18398.75
$ws.Range("L83").Value = 20350
$ws.Range("M83").Value = -13406.75
$ws.Range("N83").Value = -30334

$ws.Range("H102").Value = 2076.6667
$ws.Range("I102").Value = 1811.1666
$ws.Range("J102").Value = 2430.6667
$ws.Range("K102").Value = 1811.1666
$ws.Range("L102").Value = 2430.6667
$ws.Range("M102").Value = -189.1666
$ws.Range("N102").Value = -5674.6667

$ws.Range("H107").Value = 1301.375
$ws.Range("I107").Value = 2004.1111
$ws.Range("J107").Value = 397.85715
$ws.Range("K107").Value = 2004.1111
$ws.Range("L107").Value = 397.85715
$ws.Range("M107").Value = -84.11110000000008
$ws.Range("N107").Value = -4237.85715

$ws.Range("H128").Value = 40000
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").Value = ""

$ws.Range("H132").Value = 3719.8635
$ws.Range("I132").Value = 3142.6428
$ws.Range("K132").Value = 9427.928400000001
$ws.Range("M132").Value = -6897.928400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1564.2222
$ws.Range("I16").Value = 1630.8462
$ws.Range("J16").Value = 1391
$ws.Range("K16").Value = 1630.8462
$ws.Range("L16").Value = 1391
$ws.Range("M16").Value = -1460.8462
$ws.Range("N16").Value = -1731

$ws.Range("H54").Value = 40000
$ws.Range("J54").Value = 40000
$ws.Range("L54").Value = 40000
$ws.Range("N54").Value = -41288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = 0

$ws.Range("H12").Value = 10505.5

$ws.Range("H19").Value = 2900
$ws.Range("I19").Value = 900
$ws.Range("K19").Value = 900
$ws.Range("M19").Value = -726

$ws.Range("H20").Value = 63977.25
$ws.Range("I20").Value = 55705
$ws.Range("K20").Value = 55705
$ws.Range("M20").Value = -55465

$ws.Range("H21").Value = 32082.666
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 32082.666
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = 32082.666
$ws.Range("N21").Value = -32552.666

$ws.Range("H22").Value = 15007.5
$ws.Range("I22").Value = 10000
$ws.Range("K22").Value = 10000
$ws.Range("M22").Value = -9707

$ws.Range("H25").Value = 24027
$ws.Range("J25").Value = 24027
$ws.Range("L25").Value = 24027
$ws.Range("N25").Value = -24613

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = ""

$ws.Range("H34").Value = 46425
$ws.Range("I34").Value = 37375.332
$ws.Range("J34").Value = 59999.5
$ws.Range("K34").Value = 37375.332
$ws.Range("L34").Value = 59999.5
$ws.Range("M34").Value = -37172.332
$ws.Range("N34").Value = -60405.5

$ws.Range("H35").Value = 32082.666
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 32082.666
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = ""
$ws.Range("M35").Value = 32082.666
$ws.Range("N35").Value = -32662.666

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = ""

$ws.Range("H49").Value = 10000
$ws.Range("J49").Value = 10000
$ws.Range("L49").Value = 10000
$ws.Range("N49").Value = -10460

$ws.Range("H107").Value = 397
$ws.Range("J107").Value = 443.7143
$ws.Range("L107").Value = 1331.1429
$ws.Range("N107").Value = -5171.1429

$ws.Range("H122").Value = 3942.2307
$ws.Range("I122").Value = 4318.5
$ws.Range("J122").Value = 3340.2
$ws.Range("K122").Value = 12955.5
$ws.Range("L122").Value = 10020.6
$ws.Range("M122").Value = -10505.5
$ws.Range("N122").Value = -14920.6

$ws.Range("H123").Value = 57597.2
$ws.Range("J123").Value = 57597.2
$ws.Range("L123").Value = 57597.2
$ws.Range("N123").Value = -67397.2
